$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Category values: replace the old single "School Items, stationary" category
#     with per-product category values, adding two new product categories (H3, H4) ---
$ws.Range("H2").Value = "Electronics"
$ws.Range("H3").Value = "Toys"
$ws.Range("H4").Value = "Electronics, Toys"

# H2:H4 should carry the same Text-formatted, wrap-text style used elsewhere for
# category-like text columns.
$ws.Range("H2:H4").NumberFormat = "@"
$ws.Range("H2:H4").WrapText = $true

# --- F column (is_available) loses the stray number-format style it had and
#     goes back to the plain wrap-text style shared by most other columns ---
$ws.Range("A1").Copy()
$ws.Range("F1:F4").PasteSpecial(-4122)

# Row 2 no longer needs the taller 29pt height now that its content fits the
# default row height.
$ws.Rows.Item(2).AutoFit()

# Active selection moves to H16
$ws.Range("H16").Select()
